$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry forward the same date styling used by the existing log rows
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New finger exercise entries for lecture 11 (row 10) and Problem Set 2 (row 11)
$ws.Range("A10").Value2 = 45811
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 37
$ws.Range("F10").Value = "CS introduction lecture11"

$ws.Range("A11").Value2 = 45811
$ws.Range("B11").Value = 19
$ws.Range("C11").Value = 37
$ws.Range("F11").Value = "Problem Set 2"

# Update the selected cell to reflect where the user ended up after editing
$ws.Range("A12").Select()
